$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New partial rows (names typed in column A only) - rows 7-11
$ws.Range("A7").Value = "Noah da Mota"
$ws.Range("A8").Value = "Carlos Eduardo Santos"
$ws.Range("A9").Value = "Dr. Felipe Farias"
$ws.Range("A10").Value = "Luigi Barros"
$ws.Range("A11").Value = "Stella da Luz"

# New fully filled rows - rows 12-13
$ws.Range("A12").Value = "Noah da Mota"
$ws.Range("B12").Value = 747.91
$ws.Range("C12").Value = "472.963.815-82"
$ws.Range("D12").Value = "07/04/2024"
$ws.Range("E12").Value = "pendente"

$ws.Range("A13").Value = "Carlos Eduardo Santos"
$ws.Range("B13").Value = 123.2
$ws.Range("C13").Value = "231.659.708-40"
$ws.Range("D13").Value = "07/12/2024"
$ws.Range("E13").Value = "pendente"

# Update selection / top-left-cell view state
$ws.Range("B1").Select()
$ws.Range("C9").Select()
